$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.261.12'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.548.18'
$ws.Range('E3').Value = '  +3.24%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''569.13'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '''147.72'
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.587'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').Value = '2.546.01'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').Value = '''5.61'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = '''0.354'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '''27.48'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').Value = '3.004.79'
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('D16').Value = '63.139.87'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = '2.550.95'
$ws.Range('E18').Value = '  +3.37%  '
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = '''337.41'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''65.44'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '''1.63'
$ws.Range('E25').Value = '  +8.93%  '
$ws.Range('E26').Value = '  -2.39%  '
$ws.Range('E27').Value = '  +12.15%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''8.46'
$ws.Range('E28').Value = '  +4.61%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '''7.38'
$ws.Range('E30').Value = '  +8.26%  '
$ws.Range('D31').Value = '0.0₃0821'
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('D33').Value = '''178.54'
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('E34').Value = '  +4.00%  '
$ws.Range('D35').Value = '''411.12'
$ws.Range('E35').Value = '  +11.63%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').Value = '''19.01'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').Value = '''4.40'
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +4.34%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').Value = '''39.26'
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('D43').Value = '''153.46'
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('D45').Value = '''21.01'
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('D46').Value = '''0.604'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '''0.0522'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '''0.0238'
$ws.Range('E49').Value = '  +5.55%  '
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('E51').Value = '  +2.22%  '
